$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows above row 191, pushing the existing rows 191-203
# down to 194-206.
$ws.Rows.Item(191).Insert()
$ws.Rows.Item(191).Insert()
$ws.Rows.Item(191).Insert()

# Populate the 3 new rows (191-193) with the latest weekly price records for
# "Agrícola del Norte S.A. de Arica" - Coliflor, dated 45265 (2023-12-05).

# Row 191
$ws.Cells.Item(191, 1).Value = 1
$ws.Cells.Item(191, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(191, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(191, 4).Value = 45265
$ws.Cells.Item(191, 5).Value = 15
$ws.Cells.Item(191, 6).Value = 100112008
$ws.Cells.Item(191, 7).Value = "Coliflor"
$ws.Cells.Item(191, 8).Value = "Sin especificar"
$ws.Cells.Item(191, 9).Value = "Primera"
$ws.Cells.Item(191, 10).Value = 800
$ws.Cells.Item(191, 11).Value = 800
$ws.Cells.Item(191, 12).Value = 900
$ws.Cells.Item(191, 13).Value = 850
$ws.Cells.Item(191, 14).Value = "$/unidad"
$ws.Cells.Item(191, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(191, 16).Value = 850
$ws.Cells.Item(191, 17).Value = 1
$ws.Cells.Item(191, 18).Value = "Hortaliza"

# Row 192
$ws.Cells.Item(192, 1).Value = 1
$ws.Cells.Item(192, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(192, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(192, 4).Value = 45265
$ws.Cells.Item(192, 5).Value = 15
$ws.Cells.Item(192, 6).Value = 100112008
$ws.Cells.Item(192, 7).Value = "Coliflor"
$ws.Cells.Item(192, 8).Value = "Sin especificar"
$ws.Cells.Item(192, 9).Value = "Segunda"
$ws.Cells.Item(192, 10).Value = 1200
$ws.Cells.Item(192, 11).Value = 700
$ws.Cells.Item(192, 12).Value = 800
$ws.Cells.Item(192, 13).Value = 750
$ws.Cells.Item(192, 14).Value = "$/unidad"
$ws.Cells.Item(192, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(192, 16).Value = 750
$ws.Cells.Item(192, 17).Value = 1
$ws.Cells.Item(192, 18).Value = "Hortaliza"

# Row 193
$ws.Cells.Item(193, 1).Value = 1
$ws.Cells.Item(193, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(193, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(193, 4).Value = 45265
$ws.Cells.Item(193, 5).Value = 15
$ws.Cells.Item(193, 6).Value = 100112008
$ws.Cells.Item(193, 7).Value = "Coliflor"
$ws.Cells.Item(193, 8).Value = "Sin especificar"
$ws.Cells.Item(193, 9).Value = "Tercera"
$ws.Cells.Item(193, 10).Value = 1200
$ws.Cells.Item(193, 11).Value = 500
$ws.Cells.Item(193, 12).Value = 600
$ws.Cells.Item(193, 13).Value = 550
$ws.Cells.Item(193, 14).Value = "$/unidad"
$ws.Cells.Item(193, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(193, 16).Value = 550
$ws.Cells.Item(193, 17).Value = 1
$ws.Cells.Item(193, 18).Value = "Hortaliza"
